# Update gh-pages output (generated at 456a3b4)
# Applies:
#   - refreshed "interested" counts (column F) on several existing rows
#     in both the "展览" and "全部类型" sheets
#   - a newly scraped show ("合肥·咒术回战only") inserted into both sheets

$wb = $excel.ActiveWorkbook

function Set-LiteralText {
    # Writes $text into $cell as a literal (non-formula, non-date-parsed)
    # string value -- mirrors the inlineStr cells already in the sheet.
    param($cell, [string]$text)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

function Fix-IndexStyle {
    # Rows.Insert() leaves the brand-new row's "#" cell (column A) without
    # the bold/bordered style used throughout the table; restore it by
    # copying just the formatting from the row above.
    param($ws, [int]$row)
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------
# Sheet "展览": refresh F (interested-count) values on existing rows
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(4, 6).Value = 237
$ws1.Cells.Item(6, 6).Value = 9959
$ws1.Cells.Item(7, 6).Value = 900
$ws1.Cells.Item(10, 6).Value = 4993
$ws1.Cells.Item(13, 6).Value = 177
$ws1.Cells.Item(15, 6).Value = 62
$ws1.Cells.Item(16, 6).Value = 31
$ws1.Cells.Item(18, 6).Value = 573

# Insert a new row 21 for "合肥·咒术回战only"; old row 21 shifts to row 22.
$ws1.Rows.Item(21).Insert()
Fix-IndexStyle $ws1 21

$ws1.Cells.Item(21, 1).Value = 20
Set-LiteralText $ws1.Cells.Item(21, 2) "2024-07-28"
Set-LiteralText $ws1.Cells.Item(21, 3) "合肥·咒术回战only"
Set-LiteralText $ws1.Cells.Item(21, 4) "清河路19号 依立腾工业园区"
Set-LiteralText $ws1.Cells.Item(21, 5) "2024.07.28 09:30-07.28 17:30"
$ws1.Cells.Item(21, 6).Value = 2
$ws1.Cells.Item(21, 7).Value = 60
Set-LiteralText $ws1.Cells.Item(21, 8) "https://show.bilibili.com/platform/detail.html?id=86520"
Set-LiteralText $ws1.Cells.Item(21, 9) "//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png"

# Old row 21 ("第七届环形宇宙") now lives at row 22: bump its index + count.
$ws1.Cells.Item(22, 1).Value = 21
$ws1.Cells.Item(22, 6).Value = 1498

# ---------------------------------------------------------------------
# Sheet "全部类型": same F refreshes (rows offset by one extra entry)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(5, 6).Value = 237
$ws4.Cells.Item(7, 6).Value = 9959
$ws4.Cells.Item(8, 6).Value = 900
$ws4.Cells.Item(11, 6).Value = 4993
$ws4.Cells.Item(14, 6).Value = 177
$ws4.Cells.Item(16, 6).Value = 62
$ws4.Cells.Item(17, 6).Value = 31
$ws4.Cells.Item(19, 6).Value = 573

# Insert a new row 22 for "合肥·咒术回战only"; rows 22-23 shift to 23-24.
$ws4.Rows.Item(22).Insert()
Fix-IndexStyle $ws4 22

$ws4.Cells.Item(22, 1).Value = 21
Set-LiteralText $ws4.Cells.Item(22, 2) "2024-07-28"
Set-LiteralText $ws4.Cells.Item(22, 3) "合肥·咒术回战only"
Set-LiteralText $ws4.Cells.Item(22, 4) "清河路19号 依立腾工业园区"
Set-LiteralText $ws4.Cells.Item(22, 5) "2024.07.28 09:30-07.28 17:30"
$ws4.Cells.Item(22, 6).Value = 2
$ws4.Cells.Item(22, 7).Value = 60
Set-LiteralText $ws4.Cells.Item(22, 8) "https://show.bilibili.com/platform/detail.html?id=86520"
Set-LiteralText $ws4.Cells.Item(22, 9) "//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png"

# Old row 22 ("第七届环形宇宙") now lives at row 23: bump index + count.
$ws4.Cells.Item(23, 1).Value = 22
$ws4.Cells.Item(23, 6).Value = 1498

# Old row 23 ("首届包河留声机音乐节...") now lives at row 24: bump index only.
$ws4.Cells.Item(24, 1).Value = 23
